# Generate Report for Handoff
# Adds two new processed files (png screenshots) alongside the existing
# markdown handoff file, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

function Remove-HyperlinkAt($ws, $addr) {
    $toDelete = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $toDelete += $hl
        }
    }
    foreach ($hl in $toDelete) {
        $hl.Delete()
    }
}

$mdNew      = "1d060010-c757-44d6-aed5-405e6b3c6add.md"
$png1       = "48bd2d06-9678-489a-b5c8-a2f53de90aa4.png"
$png2       = "8a1f5cdc-175f-44ca-99a1-807ace666f66.png"

$zhcnXlf    = "1d060010-c757-44d6-aed5-405e6b3c6add.1baf98972257d5ffe3e291c2f6743e81cd03c317.zh-cn.xlf"
$dedeXlf    = "1d060010-c757-44d6-aed5-405e6b3c6add.1baf98972257d5ffe3e291c2f6743e81cd03c317.de-de.xlf"

$png1Target = "059e19a1971e662d4b7f4737b3c09e874ed09852.png"
$png2Target = "ce661e7b01eaf9c1b812fe77eba05d28f888e198.png"

$overviewDate  = "2016-48-17 22:48:41"
$zhcnHandoffDt = "2016-03-17 22:48:38"
$dedeHandoffDt = "2016-03-17 22:48:41"

$epoch        = "0001-01-01 00:00:00"
$readyStatus  = "Ready for handoff"
$includeRsn   = "Include"
$isDepRsn     = "IsDependency"
$depFrom      = "e2e\" + $mdNew

$srcBase   = "https://github.com/OpenLocalizationTest/oltest/blob/cadb00bdf4058bafce939b9af00f9a5207bcf73e/e2e/"
$zhBase    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5a9d5ac742aa63c9b0caeaaee58500f08511e186/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deBase    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d695193064d8fb305621a3c6b54e7718b2e25e23/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Remove-HyperlinkAt $wsOverview '$A$2'
$a2 = $wsOverview.Range("A2")
$a2.Value = $mdNew
$wsOverview.Hyperlinks.Add($a2, $srcBase + $mdNew, $null, $null, $mdNew)

$wsOverview.Range("B2").Value = $readyStatus
$wsOverview.Range("C2").Value = $readyStatus
$wsOverview.Range("D2").Value = $overviewDate

$a3 = $wsOverview.Range("A3")
$a3.Value = $png1
$wsOverview.Hyperlinks.Add($a3, $srcBase + $png1, $null, $null, $png1)
$wsOverview.Range("B3").Value = $readyStatus
$wsOverview.Range("C3").Value = $readyStatus
$wsOverview.Range("D3").Value = $overviewDate

$a4 = $wsOverview.Range("A4")
$a4.Value = $png2
$wsOverview.Hyperlinks.Add($a4, $srcBase + $png2, $null, $null, $png2)
$wsOverview.Range("B4").Value = $readyStatus
$wsOverview.Range("C4").Value = $readyStatus
$wsOverview.Range("D4").Value = $overviewDate

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Remove-HyperlinkAt $wsZhCn '$A$2'
Remove-HyperlinkAt $wsZhCn '$B$2'
Remove-HyperlinkAt $wsZhCn '$D$2'

$a2 = $wsZhCn.Range("A2")
$a2.Value = $mdNew
$wsZhCn.Hyperlinks.Add($a2, $srcBase + $mdNew, $null, $null, $mdNew)

$b2 = $wsZhCn.Range("B2")
$b2.Value = ".md"
$wsZhCn.Hyperlinks.Add($b2, $srcBase + $mdNew, $null, $null, ".md")

$wsZhCn.Range("C2").Value = $readyStatus

$d2 = $wsZhCn.Range("D2")
$d2.Value = $zhcnXlf
$wsZhCn.Hyperlinks.Add($d2, $zhBase + $zhcnXlf, $null, $null, $zhcnXlf)

$wsZhCn.Range("E2").Value = $zhcnHandoffDt
$wsZhCn.Range("H2").Value = $epoch
$wsZhCn.Range("I2").Value = $includeRsn

$a3 = $wsZhCn.Range("A3")
$a3.Value = $png1
$wsZhCn.Hyperlinks.Add($a3, $srcBase + $png1, $null, $null, $png1)

$b3 = $wsZhCn.Range("B3")
$b3.Value = ".png"
$wsZhCn.Hyperlinks.Add($b3, $srcBase + $png1, $null, $null, ".png")

$wsZhCn.Range("C3").Value = $readyStatus

$d3 = $wsZhCn.Range("D3")
$d3.Value = $png1Target
$wsZhCn.Hyperlinks.Add($d3, $zhBase + $png1Target, $null, $null, $png1Target)

$wsZhCn.Range("E3").Value = $zhcnHandoffDt
$wsZhCn.Range("H3").Value = $epoch
$wsZhCn.Range("I3").Value = $isDepRsn
$wsZhCn.Range("J3").Value = $depFrom

$a4 = $wsZhCn.Range("A4")
$a4.Value = $png2
$wsZhCn.Hyperlinks.Add($a4, $srcBase + $png2, $null, $null, $png2)

$b4 = $wsZhCn.Range("B4")
$b4.Value = ".png"
$wsZhCn.Hyperlinks.Add($b4, $srcBase + $png2, $null, $null, ".png")

$wsZhCn.Range("C4").Value = $readyStatus

$d4 = $wsZhCn.Range("D4")
$d4.Value = $png2Target
$wsZhCn.Hyperlinks.Add($d4, $zhBase + $png2Target, $null, $null, $png2Target)

$wsZhCn.Range("E4").Value = $zhcnHandoffDt
$wsZhCn.Range("H4").Value = $epoch
$wsZhCn.Range("I4").Value = $isDepRsn
$wsZhCn.Range("J4").Value = $depFrom

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Remove-HyperlinkAt $wsDeDe '$A$2'
Remove-HyperlinkAt $wsDeDe '$B$2'
Remove-HyperlinkAt $wsDeDe '$D$2'

$a2 = $wsDeDe.Range("A2")
$a2.Value = $mdNew
$wsDeDe.Hyperlinks.Add($a2, $srcBase + $mdNew, $null, $null, $mdNew)

$b2 = $wsDeDe.Range("B2")
$b2.Value = ".md"
$wsDeDe.Hyperlinks.Add($b2, $srcBase + $mdNew, $null, $null, ".md")

$wsDeDe.Range("C2").Value = $readyStatus

$d2 = $wsDeDe.Range("D2")
$d2.Value = $dedeXlf
$wsDeDe.Hyperlinks.Add($d2, $deBase + $dedeXlf, $null, $null, $dedeXlf)

$wsDeDe.Range("E2").Value = $dedeHandoffDt
$wsDeDe.Range("H2").Value = $epoch
$wsDeDe.Range("I2").Value = $includeRsn

$a3 = $wsDeDe.Range("A3")
$a3.Value = $png1
$wsDeDe.Hyperlinks.Add($a3, $srcBase + $png1, $null, $null, $png1)

$b3 = $wsDeDe.Range("B3")
$b3.Value = ".png"
$wsDeDe.Hyperlinks.Add($b3, $srcBase + $png1, $null, $null, ".png")

$wsDeDe.Range("C3").Value = $readyStatus

$d3 = $wsDeDe.Range("D3")
$d3.Value = $png1Target
$wsDeDe.Hyperlinks.Add($d3, $deBase + $png1Target, $null, $null, $png1Target)

$wsDeDe.Range("E3").Value = $dedeHandoffDt
$wsDeDe.Range("H3").Value = $epoch
$wsDeDe.Range("I3").Value = $isDepRsn
$wsDeDe.Range("J3").Value = $depFrom

$a4 = $wsDeDe.Range("A4")
$a4.Value = $png2
$wsDeDe.Hyperlinks.Add($a4, $srcBase + $png2, $null, $null, $png2)

$b4 = $wsDeDe.Range("B4")
$b4.Value = ".png"
$wsDeDe.Hyperlinks.Add($b4, $srcBase + $png2, $null, $null, ".png")

$wsDeDe.Range("C4").Value = $readyStatus

$d4 = $wsDeDe.Range("D4")
$d4.Value = $png2Target
$wsDeDe.Hyperlinks.Add($d4, $deBase + $png2Target, $null, $null, $png2Target)

$wsDeDe.Range("E4").Value = $dedeHandoffDt
$wsDeDe.Range("H4").Value = $epoch
$wsDeDe.Range("I4").Value = $isDepRsn
$wsDeDe.Range("J4").Value = $depFrom

Write-Host "Report generated for handoff"
